$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.8871922249200372, 0.5969110571806628, 1.3619404499586376, 2.669885617651782)"
$ws.Range("C2").Value = "NIG(1.407552908362087, 1.0940782112350167, 3.140343901205358, 5.214953775489716)"
$ws.Range("D2").Value = "NCT(2.2708972336861946, 1.5065209268749937, -0.8301230842375278, 2.2836820928583794)"
$ws.Range("E2").Value = "JSU(-1.4099427221920178, 1.19363204025625, 2.161080550921044, 4.480679254835424)"
